$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the data table -------------------------------------------------
# Old layout (A:E) = 구역 이름 | 사고 위험도 | 사고 규모 | 최소 필요 인원 | 최대 필요 인원
#   with 6 data rows (rows 2-7)
# New layout (A:E) = region name | min required | max required | lower bound (demand) | accident scale
#   with only 2 data rows (rows 2-3), regions renamed region0/region1/region2

$ws.Range("A1:E7").ClearContents()
$ws.Rows("4:7").Delete()

$ws.Range("A1").Value = "region0"
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 50

$ws.Range("A2").Value = "region1"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 30

$ws.Range("A3").Value = "region2"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 50

# --- View / selection changes -----------------------------------------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("E7").Select()
